# Add a new "Question 4" sheet (geography query for the map) after "Suggestion4".
$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $afterSheet)
$ws.Name = "Question 4"

# Title / question row (A1) - gets its own "apply alignment" style once wrap
# text is explicitly toggled, matching the rest of the workbook's "question"
# header cells. Re-run AutoFit afterwards so the embedded line break in the
# text doesn't leave a custom row height behind.
$ws.Range("A1").Value = "Question 4: Where should the customer accquisition funds should be spent to increase`r`nthe rentals?"
$ws.Range("A1").WrapText = $false
$ws.Rows.Item(1).AutoFit() | Out-Null

# Column headers
$ws.Range("A2").Value = "country"
$ws.Range("B2").Value = "rentals"

# Country / rentals data pulled from the SQL geography query results
$data = @"
India	1572
China	1426
United States	968
Japan	825
Mexico	796
Brazil	748
Russian Federation	713
Philippines	568
Turkey	388
Indonesia	367
Argentina	352
Nigeria	352
Taiwan	305
South Africa	285
Iran	225
United Kingdom	219
Poland	203
Germany	196
Italy	189
Vietnam	172
Venezuela	170
Egypt	161
Colombia	159
Ukraine	158
Spain	142
Canada	137
South Korea	135
Netherlands	134
Pakistan	128
Saudi Arabia	121
Yemen	117
Peru	112
Israel	99
France	96
Thailand	96
Bangladesh	95
Algeria	90
Ecuador	87
Malaysia	85
United Arab Emirates	84
Tanzania	83
Mozambique	80
Austria	78
Dominican Republic	77
Chile	71
Morocco	71
Belarus	66
Paraguay	62
Romania	62
Puerto Rico	61
Switzerland	61
Latvia	60
Azerbaijan	57
Yugoslavia	57
Cameroon	54
French Polynesia	54
Greece	54
Kenya	54
Sudan	54
Cambodia	53
Angola	52
Kazakstan	52
Myanmar	52
Bulgaria	50
Congo, The Democratic Republic of the	50
Oman	50
Bolivia	47
Runion	46
Anguilla	35
Brunei	35
Greenland	34
Holy See (Vatican City State)	34
Moldova	34
Sweden	34
Zambia	33
Chad	32
Virgin Islands, U.S.	32
Nauru	31
North Korea	31
Estonia	30
Gambia	30
Hong Kong	30
Sri Lanka	30
Czech Republic	29
Hungary	29
Faroe Islands	28
Liechtenstein	28
Malawi	28
Iraq	27
Turkmenistan	27
Finland	26
Kuwait	26
Slovakia	26
Tuvalu	26
Armenia	25
Bahrain	25
Saint Vincent and the Grenadines	25
Senegal	25
Lithuania	24
New Zealand	24
Ethiopia	23
Tunisia	23
French Guiana	22
Madagascar	22
Nepal	22
American Samoa	20
Afghanistan	18
Tonga	18
"@

$rows = $data -split "`n"
$r = 3
foreach ($line in $rows) {
    $parts = $line -split "`t"
    $ws.Cells.Item($r, 1).Value = $parts[0]
    $ws.Cells.Item($r, 2).Value = [int]$parts[1]
    $r = $r + 1
}

# Match the author's final selection before saving.
$ws.Range("F7").Select() | Out-Null
